# Apply the "Automatic update of files" edit.
#
# The commit cyclically rotates a fixed set of columns across data rows 4-11
# (8 observation rows) by one position: each row's data for these columns is
# replaced by the data that used to live in the row above it, and row 4
# wraps around to take what used to be in row 11.
#
# Columns involved in the rotation:
#   A, B, E, F, G, H, P, Q, R, Z, AB, AC, AM, AO
#
# All other columns / cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> column index map for the columns that rotate.
$colMap = @{
    "A"  = 1
    "B"  = 2
    "E"  = 5
    "F"  = 6
    "G"  = 7
    "H"  = 8
    "P"  = 16
    "Q"  = 17
    "R"  = 18
    "Z"  = 26
    "AB" = 28
    "AC" = 29
    "AM" = 39
    "AO" = 41
}

$rows = @(4, 5, 6, 7, 8, 9, 10, 11)

# 1) Snapshot the original values for every rotating cell BEFORE making any
#    changes, so later writes don't clobber values we still need to read.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $colMap.Keys) {
        $idx = $colMap[$col]
        $rowData[$col] = $ws.Cells.Item($r, $idx).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write each row's new values, pulling from the row above it (cyclically,
#    row 4 pulls from row 11).
for ($i = 0; $i -lt $rows.Count; $i++) {
    $destRow = $rows[$i]
    if ($i -eq 0) {
        $srcRow = $rows[$rows.Count - 1]
    } else {
        $srcRow = $rows[$i - 1]
    }

    $srcData = $snapshot[$srcRow]

    foreach ($col in $colMap.Keys) {
        $idx = $colMap[$col]
        $val = $srcData[$col]
        if ($val -eq $null) {
            $ws.Cells.Item($destRow, $idx).ClearContents()
        } else {
            $ws.Cells.Item($destRow, $idx).Value = $val
        }
    }
}
